$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "66.884.13"
Set-TextValue "D3" "2.612.01"
Set-TextValue "E3" "  -1.55%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "586.50"
Set-TextValue "E5" "  -1.62%  "
Set-TextValue "D6" "165.38"
Set-TextValue "E6" "  -1.44%  "
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "D8" "0.528"
Set-TextValue "E8" "  -3.69%  "
Set-TextValue "D9" "2.610.75"
Set-TextValue "E9" "  -1.60%  "
Set-TextValue "E10" "  -3.64%  "
Set-TextValue "E11" "  +0.41%  "
Set-TextValue "D12" "0.365"
Set-TextValue "E12" "  -0.18%  "
Set-TextValue "D13" "5.18"
Set-TextValue "E13" "  -1.79%  "
Set-TextValue "D14" "27.18"
Set-TextValue "E14" "  -3.31%  "
Set-TextValue "D15" "3.089.79"
Set-TextValue "D16" "0.0000179"
Set-TextValue "E16" "  -3.04%  "
Set-TextValue "D17" "66.767.74"
Set-TextValue "E17" "  -1.17%  "
Set-TextValue "D18" "2.586.38"
Set-TextValue "E18" "  -2.53%  "
Set-TextValue "D19" "11.65"
Set-TextValue "E19" "  -4.14%  "
Set-TextValue "D20" "7.78"
Set-TextValue "E20" "  -4.91%  "
Set-TextValue "D21" "354.49"
Set-TextValue "E21" "  -2.48%  "
Set-TextValue "D22" "4.26"
Set-TextValue "E22" "  -3.38%  "
Set-TextValue "D23" "4.62"
Set-TextValue "E23" "  -3.90%  "
Set-TextValue "D24" "10.48"
Set-TextValue "E24" "  -5.13%  "
Set-TextValue "E25" "  -0.05%  "
Set-TextValue "E26" "  -5.68%  "
Set-TextValue "D27" "69.30"
Set-TextValue "E27" "  -2.80%  "
Set-TextValue "D28" "2.745.97"
Set-TextValue "E28" "  -1.62%  "
Set-TextValue "E29" "  +0.05%  "
Set-TextValue "D30" "0.0₃0991"
Set-TextValue "E30" "  -3.75%  "
Set-TextValue "D31" "539.43"
Set-TextValue "E31" "  -3.40%  "
Set-TextValue "D32" "8.15"
Set-TextValue "E32" "  +1.09%  "
Set-TextValue "D33" "1.33"
Set-TextValue "E33" "  -4.55%  "
Set-TextValue "E34" "  -3.36%  "
Set-TextValue "E35" "  -1.71%  "
Set-TextValue "D36" "1.00"
Set-TextValue "E37" "  -5.27%  "
Set-TextValue "D38" "158.31"
Set-TextValue "E38" "  +0.99%  "
Set-TextValue "D39" "18.87"
Set-TextValue "E39" "  -2.66%  "
Set-TextValue "D40" "0.363"
Set-TextValue "E40" "  -2.67%  "
Set-TextValue "E41" "  +1.84%  "
Set-TextValue "D42" "1.79"
Set-TextValue "E42" "  -1.89%  "
Set-TextValue "D43" "5.11"
Set-TextValue "E43" "  -3.97%  "
Set-TextValue "D45" "2.40"
Set-TextValue "E45" "  -5.79%  "
Set-TextValue "D46" "0.0₆0293"
Set-TextValue "E46" "  -2.12%  "
Set-TextValue "B47" "Aave"
Set-TextValue "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D47" "150.51"
Set-TextValue "E47" "  -2.60%  "
Set-TextValue "B48" "ARBITRUM"
Set-TextValue "C48" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D48" "0.574"
Set-TextValue "E48" "  -3.67%  "
Set-TextValue "D49" "3.74"
Set-TextValue "E49" "  -3.68%  "
Set-TextValue "D50" "1.70"
Set-TextValue "E50" "  -2.06%  "
Set-TextValue "E51" "  -1.51%  "
